$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Drop the stale "_GoBack" bookmark that currently sits in the
#    "Date : 02/03/2020" heading (it marks Word's last-edit spot and
#    needs to move to wherever we edit next).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Find the "Séance N°X" heading paragraph and turn the trailing
#    "X" into its own run reading "5" (-> "Séance N°" + "5").
# ------------------------------------------------------------------
$seance = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Séance N°X*") {
        $seance = $p
        break
    }
}

$pEnd = $seance.Range.End

# The paragraph mark itself occupies (pEnd-1, pEnd); the last visible
# character ("X") sits right before it, at (pEnd-2, pEnd-1).
$xRange = $d.Range($pEnd - 2, $pEnd - 1)
$xRange.Delete()

$insertPos = $pEnd - 2
$d.Range($insertPos, $insertPos).InsertAfter("5")

# ------------------------------------------------------------------
# 3. Re-create "_GoBack" right after the new "5", at the very end of
#    the "Séance N°5" paragraph (an empty bookmark, same as Word
#    leaves behind after typing at that spot).
#
#    Inserting a collapsed bookmark exactly at a paragraph's final
#    position (immediately before its paragraph mark) lands on the
#    very start of the document, so a short-lived placeholder run is
#    used to push that boundary out of the way, the bookmark is
#    planted in front of it, and the placeholder is removed again.
# ------------------------------------------------------------------
$seance = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Séance N°5*") {
        $seance = $p
        break
    }
}

$endPos = $seance.Range.End - 1
$d.Range($endPos, $endPos).InsertAfter("ZZZ")
$d.Bookmarks.Add("_GoBack", $d.Range($endPos, $endPos))
$d.Range($endPos, $endPos + 3).Delete()
